# edit.ps1 -- apply the "raven.docx" revision:
#   1. Append "  (This is a change – Version for main branch)" to the end
#      of the first paragraph, in three red-colored runs (matching the
#      author's incremental-typing run boundaries), after widening the
#      original run's text to keep its two trailing spaces.
#   2. Delete the trailing "ank God almighty, we are free at last."
#      paragraph (the start of that sentence, "Th", already lives in the
#      previous run - this whole paragraph goes away).
#   3. Drop the styles that became unused/unreferenced once that
#      paragraph (and the page's heading/list scaffolding) went away.

$d = $word.ActiveDocument

# --- 1. First paragraph: widen the base run, then add 3 red runs ------------
$r = $d.Paragraphs.Item(1).Range
$null = $r.MoveEnd(1, -1)              # exclude the paragraph mark itself

$r.InsertAfter("  ")
$r.Collapse(0)

$r.InsertAfter([char]0x0028 + "This is a change " + [char]0x2013 + " Ve")
$r.Font.Color = 255                    # wdColorRed / RGB(255,0,0) -> FF0000
$r.Collapse(0)

$r.InsertAfter("rsion for main branch")
$r.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter([char]0x0029)
$r.Font.Color = 255
$r.Collapse(0)

# --- 2. Remove the last paragraph ("ank God almighty, we are free at last.") -
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Delete()

# --- 3. Remove styles that are now unused ------------------------------------
$obsoleteStyleNames = @(
    "Heading 2",
    "Heading 4",
    "apple-converted-space",
    "Hyperlink",
    "Heading 2 Char",
    "Heading 4 Char",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)

$count = $d.Styles.Count
$toDelete = @()
for ($i = 1; $i -le $count; $i++) {
    $nm = $d.Styles.Item($i).NameLocal
    if ($obsoleteStyleNames -contains $nm) {
        $toDelete += $i
    }
}

# Delete from the highest index down so earlier indices stay valid.
$toDelete = $toDelete | Sort-Object -Descending
foreach ($i in $toDelete) {
    $d.Styles.Item($i).Delete()
}

Write-Output "done"
